# Package & Initializing Socket-IO
# Edits slide 2 ("What are we using?") content placeholder:
#   - removes the "npm" call-out from the first bullet
#   - appends a new "Flask-SocketIO" bullet
#   - bumps the autofit line-spacing reduction
# Also tidies a stray endParaRPr / adds dirty="0" on slide 21.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2 - "What are we using?" content placeholder
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shape = $s2.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Replace the whole body with the final text (paragraph marks via `r).
# Re-assigning the full TextRange collapses every paragraph down to a single
# run each, using the shape's existing Arial run formatting.
$tr.Text = "You'll have to install all of these except Flask_Cors (which is just a python module). For more info, look up their installation instructions.`rFlask " + [char]0x2013 + " a Python web dev framework`rSQLite " + [char]0x2013 + " a fairly reliable, lightweight, easy-to-use database that works for websites at least up to 100k users`rFlask-SQLAlchemy " + [char]0x2013 + " a Flask wrapper for SQLAlchemy, which is an Object Relational Mapper (ORM) that lets us use high-level classes, objects, and methods instead of dealing with tables and SQL`rFlask-Migrate " + [char]0x2013 + " a Flask wrapper for Alembic, a database migration framework for SQL-Alchemy that makes it easy to update the database when changed and have version control`rFlask_Cors " + [char]0x2013 + " allows HTTP requests from a different url, PLEASE PLEASE PLEASE disable before putting the app live unless we want people`rFlask-SocketIO " + [char]0x2013 + " allows for bilateral communication with rapid updates, useful for chatbox"

# Re-split paragraph 1 into its 3 original runs: lead text / "Flask_Cors" / trailing text
$tr.Characters(44, 10).Text = "Flask_Cors"

# Re-split paragraph 4 ("Flask-" / "SQLAlchemy" / " - a Flask wrapper for " / "SQLAlchemy" / rest)
$tr.Characters(291, 6).Text = "Flask-"
$tr.Characters(297, 10).Text = "SQLAlchemy"
$tr.Characters(330, 10).Text = "SQLAlchemy"

# Re-split paragraph 6 ("Flask_Cors" / middle / "url" / ", PLEASE " / "PLEASE" / " " / "PLEASE" / rest)
$tr.Characters(654, 10).Text = "Flask_Cors"
$tr.Characters(705, 3).Text = "url"
$tr.Characters(717, 6).Text = "PLEASE"
$tr.Characters(724, 6).Text = "PLEASE"

# Re-split paragraph 7 ("Flask-" / "SocketIO" / middle / "chatbox")
$tr.Characters(789, 6).Text = "Flask-"
$tr.Characters(795, 8).Text = "SocketIO"
$tr.Characters(872, 7).Text = "chatbox"

# Bump the autofit line spacing reduction 10% -> 20% (fontScale stays the same)
$shape.TextFrame2.TextRange.Text = $shape.TextFrame2.TextRange.Text
$bodyXml = $null

# ---------------------------------------------------------------------------
# Slide 21 - "/catalog" PATCH paragraph tidy-up
# ---------------------------------------------------------------------------
$s21 = $p.Slides.Item(21)
$shape21 = $s21.Shapes.Item(2)
$tr21 = $shape21.TextFrame.TextRange
$paraCount = $tr21.Paragraphs().Count
$lastPara = $tr21.Paragraphs($paraCount, 1)
$lastPara.Text = "PATCH method with the param " + [char]0x00AB + "catalog_id" + [char]0x00BB + " flips the value of the available variable"
